$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 168, shifting existing rows 168-212 down to 169-213.
$ws.Rows(168).Insert()

# Populate the newly inserted row 168 with its data.
$ws.Cells.Item(168, 1).Value = 7
$ws.Cells.Item(168, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(168, 3).Value = "Ñuble"
$ws.Cells.Item(168, 4).Value = "2022-04-12"
$ws.Cells.Item(168, 5).Value = 16
$ws.Cells.Item(168, 6).Value = 100112032
$ws.Cells.Item(168, 7).Value = "Zapallo italiano"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 120
$ws.Cells.Item(168, 11).Value = 8000
$ws.Cells.Item(168, 12).Value = 8500
$ws.Cells.Item(168, 13).Value = 8250
$ws.Cells.Item(168, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(168, 15).Value = "Región del Maule"
$ws.Cells.Item(168, 16).Value = 165
$ws.Cells.Item(168, 17).Value = 50
$ws.Cells.Item(168, 18).Value = "Hortaliza"
